{"js": "// The author fixed a capitalization typo: \"erscheint Sie unter\" should be\n// \"erscheint sie unter\" (the reflexive/personal pronoun \"sie\" referring back\n// to \"die Datei\" must be lowercase here, it is not the polite \"Sie\").\n// Find the unique sentence fragment and replace the capitalized \"Sie\" with\n// lowercase \"sie\", leaving the rest of the sentence untouched.\n\nconst body = context.document.body;\n\nconst results = body.search(\"erscheint Sie unter\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const found = results.items[0];\n  found.insertText(\"erscheint sie unter\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# The author fixed a capitalization typo: \"erscheint Sie unter\" should be\n# \"erscheint sie unter\" (the pronoun \"sie\" refers back to \"die Datei\" and\n# must be lowercase here - it is not the polite/formal \"Sie\").\n# \"Sie\"/\"sie\" occurs many times in this document (e.g. \"Klicken Sie auf die\n# Datei.\"), so we first locate the unique sentence fragment, then replace\n# only the single capital \"S\" that sits inside it with a lowercase \"s\".\n\n$d = $word.ActiveDocument\n\n$scope = $d.Content\n$find = $scope.Find\n$find.ClearFormatting()\n$find.Text = \"erscheint Sie unter\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Forward = $true\n$find.Wrap = 0\n\nif ($find.Execute()) {\n  # $scope now exactly spans \"erscheint Sie unter\"; find the \"S\" inside it\n  # and narrow to that single character so nothing else is touched.\n  $text = $scope.Text\n  $offset = $text.IndexOf(\"Sie\")\n  if ($offset -ge 0) {\n    $capStart = $scope.Start + $offset\n    $capRange = $d.Range($capStart, $capStart + 1)\n    if ($capRange.Text -ceq \"S\") {\n      $capRange.Text = \"s\"\n    }\n  }\n}\n"}
